$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feeders")

# --- Row 3 ("prosthetic" feeder) gains a real event-class mapping ---
# C3: was blank/"none" with no special format -> now "Prosthesis", formatted
#     like the other event_class_code cells (shaded fill + border) but with
#     an explicit text number format.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "Prosthesis"

# D3: new event_class_seq value, same look as the other seq cells (D2/D5).
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 6.2

# E3: feeder_description changes, now with a highlighted box around it
#     (shaded fill + border on right/top/bottom only).
$ws.Range("C2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Borders(7).LineStyle = -4142
$ws.Range("E3").Borders(10).Color = 9359529
$ws.Range("E3").Value = "Prosthetic Expenses"

# --- Row 4 ("radiology" feeder): event_class_seq 9.1 -> 7.1 ---
$ws.Range("D4").Value = 7.1

# Clear the clipboard marching ants / mirror the saved selection state.
$excel.CutCopyMode = 0
$ws.Range("E13").Select()
